$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of this producer's weekly price-history
# block (existing rows 84-91 shift down to 86-93), then populate the two
# new rows with this week's data.
$ws.Rows("84:85").Insert()

# Row 84 (newest entry, Calidad = Primera)
$ws.Range("A84").Value = 1
$ws.Range("B84").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C84").Value = "Arica y Parinacota"
$ws.Range("D84").Value = 44617
$ws.Range("E84").Value = 15
$ws.Range("F84").Value = 100112036
$ws.Range("G84").Value = "Caigua"
$ws.Range("H84").Value = "Sin especificar"
$ws.Range("I84").Value = "Primera"
$ws.Range("J84").Value = 120
$ws.Range("K84").Value = 12000
$ws.Range("L84").Value = 13000
$ws.Range("M84").Value = 12500
$ws.Range("N84").Value = "`$/caja 20 kilos"
$ws.Range("O84").Value = "Región de Arica y Parinacota"
$ws.Range("P84").Value = 625
$ws.Range("Q84").Value = 20
$ws.Range("R84").Value = "Hortaliza"

# Row 85 (newest entry, Calidad = Segunda)
$ws.Range("A85").Value = 1
$ws.Range("B85").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C85").Value = "Arica y Parinacota"
$ws.Range("D85").Value = 44617
$ws.Range("E85").Value = 15
$ws.Range("F85").Value = 100112036
$ws.Range("G85").Value = "Caigua"
$ws.Range("H85").Value = "Sin especificar"
$ws.Range("I85").Value = "Segunda"
$ws.Range("J85").Value = 130
$ws.Range("K85").Value = 9000
$ws.Range("L85").Value = 10000
$ws.Range("M85").Value = 9500
$ws.Range("N85").Value = "`$/caja 20 kilos"
$ws.Range("O85").Value = "Región de Arica y Parinacota"
$ws.Range("P85").Value = 475
$ws.Range("Q85").Value = 20
$ws.Range("R85").Value = "Hortaliza"
